$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update category image links (column E) to the new cat_images paths
$ws.Range("E4").Value  = "cat_images/shorts_apparel_merchandise.jpg"
$ws.Range("E5").Value  = "cat_images/tshirts_apparel_merchandise.jpg"
$ws.Range("E6").Value  = "cat_images/sweatshirts_apparel_merchandise.jpgs"
$ws.Range("E7").Value  = "cat_images/jerseys_apparel_merchandise.jpg"
$ws.Range("E8").Value  = "cat_images/pants_apparel_merchandise.jpg"
$ws.Range("E9").Value  = "cat_images/hats_apparel_merchandise.jpg"
$ws.Range("E10").Value = "cat_images/novelties_accessories.jpgs"
$ws.Range("E11").Value = "cat_images/mens_apparel_merchandise.jpg"
$ws.Range("E12").Value = "cat_images/womens_apparel_merchandise.jpg"
$ws.Range("E13").Value = "cat_images/kids_apparel_merchandise.jpg"
$ws.Range("E14").Value = "cat_images/youth_apparel_merchandise.jpg"
$ws.Range("E15").Value = "cat_images/toddler_apparel_merchandise.jpg"
$ws.Range("E16").Value = "cat_images/infant_apparel_merchandise.jpg"

# Re-point category names (column A) that previously used the duplicate
# lowercase/alternate string set, onto the canonical display strings
$ws.Range("A6").Value  = "Sweatshirts"
$ws.Range("A10").Value = "Novelties and Accessories"
$ws.Range("A12").Value = "Women's"

# Update the view state: scroll back to show column A (drop topLeftCell)
# and move the active selection to D19
$ws.Range("D19").Select()
